# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 56.142857
$ws.Cells.Item(11, 9).Value = 56.142857
$ws.Cells.Item(11, 11).Value = 56.142857
$ws.Cells.Item(11, 13).Value = 83.85714300000001
$ws.Cells.Item(40, 8).Value = 4890.68
$ws.Cells.Item(40, 9).Value = 2642.625
$ws.Cells.Item(40, 11).Value = 2642.625
$ws.Cells.Item(40, 13).Value = -2467.625
$ws.Cells.Item(41, 8).Value = 20001668
$ws.Cells.Item(41, 9).Value = 475
$ws.Cells.Item(41, 10).Value = 45457732
$ws.Cells.Item(41, 11).Value = 475
$ws.Cells.Item(41, 12).Value = 45457732
$ws.Cells.Item(41, 13).Value = -35
$ws.Cells.Item(41, 14).Value = -45458612
$ws.Cells.Item(62, 8).Value = 12505548
$ws.Cells.Item(62, 9).Value = 25003096
$ws.Cells.Item(62, 11).Value = 25003096
$ws.Cells.Item(62, 13).Value = -25002472
$ws.Cells.Item(65, 8).Value = 12505548
$ws.Cells.Item(65, 9).Value = 25003096
$ws.Cells.Item(65, 11).Value = 125015480
$ws.Cells.Item(65, 13).Value = -125012360
$ws.Cells.Item(69, 8).Value = 22500
$ws.Cells.Item(69, 9).Value = 20000
$ws.Cells.Item(69, 11).Value = 60000
$ws.Cells.Item(69, 13).Value = -59126
$ws.Cells.Item(72, 8).Value = 22500
$ws.Cells.Item(72, 9).Value = 20000
$ws.Cells.Item(72, 11).Value = 180000
$ws.Cells.Item(72, 13).Value = -175632
$ws.Cells.Item(96, 9).Value = 224.75
$ws.Cells.Item(96, 10).Value = 522.8570999999999
$ws.Cells.Item(96, 11).Value = 674.25
$ws.Cells.Item(96, 12).Value = 1568.5713
$ws.Cells.Item(96, 13).Value = 698.75
$ws.Cells.Item(96, 14).Value = -4314.5713
$ws.Cells.Item(113, 8).Value = 9777.6
$ws.Cells.Item(113, 9).Value = 11200.429
$ws.Cells.Item(113, 10).Value = 9011.462
$ws.Cells.Item(113, 11).Value = 11200.429
$ws.Cells.Item(113, 12).Value = 9011.462
$ws.Cells.Item(113, 13).Value = -7946.429
$ws.Cells.Item(113, 14).Value = -15519.462
$ws.Cells.Item(132, 8).Value = 6268.905
$ws.Cells.Item(132, 9).Value = 6282.35
$ws.Cells.Item(132, 11).Value = 18847.05
$ws.Cells.Item(132, 13).Value = -16317.05

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 133471.55
$ws.Cells.Item(110, 9).Value = 174195.4
$ws.Cells.Item(110, 10).Value = 2250.2222
$ws.Cells.Item(110, 11).Value = 174195.4
$ws.Cells.Item(110, 12).Value = 2250.2222
$ws.Cells.Item(110, 13).Value = -172150.4
$ws.Cells.Item(110, 14).Value = -6340.2222
$ws.Cells.Item(122, 8).Value = 4635.5
$ws.Cells.Item(122, 9).Value = 4105.9565
$ws.Cells.Item(122, 11).Value = 12317.8695
$ws.Cells.Item(122, 13).Value = -9867.869500000001
$ws.Cells.Item(132, 8).Value = 8499.25
$ws.Cells.Item(132, 9).Value = 8498
$ws.Cells.Item(132, 10).Value = 8499.666999999999
$ws.Cells.Item(132, 11).Value = 25494
$ws.Cells.Item(132, 12).Value = 25499.001
$ws.Cells.Item(132, 13).Value = -22964
$ws.Cells.Item(132, 14).Value = -30559.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 58825068
$ws.Cells.Item(20, 9).Value = 90910660
$ws.Cells.Item(20, 10).Value = 1492.8334
$ws.Cells.Item(20, 11).Value = 90910660
$ws.Cells.Item(20, 12).Value = 1492.8334
$ws.Cells.Item(20, 13).Value = -90910413
$ws.Cells.Item(20, 14).Value = -1986.8334
$ws.Cells.Item(64, 8).Value = 950.7
$ws.Cells.Item(64, 9).Value = 678
$ws.Cells.Item(64, 10).Value = 1132.5
$ws.Cells.Item(64, 11).Value = 678
$ws.Cells.Item(64, 12).Value = 1132.5
$ws.Cells.Item(64, 13).Value = -453
$ws.Cells.Item(64, 14).Value = -1582.5
$ws.Cells.Item(67, 8).Value = 950.7
$ws.Cells.Item(67, 9).Value = 678
$ws.Cells.Item(67, 10).Value = 1132.5
$ws.Cells.Item(67, 11).Value = 678
$ws.Cells.Item(67, 12).Value = 1132.5
$ws.Cells.Item(67, 13).Value = 102
$ws.Cells.Item(67, 14).Value = -2692.5
$ws.Cells.Item(134, 8).Value = 51127.914
$ws.Cells.Item(134, 10).Value = 117766.445
$ws.Cells.Item(134, 12).Value = 353299.335
$ws.Cells.Item(134, 14).Value = -358369.335

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3674.889
$ws.Cells.Item(31, 9).Value = 1634.25
$ws.Cells.Item(31, 11).Value = 1634.25
$ws.Cells.Item(31, 13).Value = -1339.25
$ws.Cells.Item(34, 8).Value = 3674.889
$ws.Cells.Item(34, 9).Value = 1634.25
$ws.Cells.Item(34, 11).Value = 1634.25
$ws.Cells.Item(34, 13).Value = -1432.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 790.4737
$ws.Cells.Item(98, 9).Value = 1091
$ws.Cells.Item(98, 10).Value = 520
$ws.Cells.Item(98, 11).Value = 3273
$ws.Cells.Item(98, 12).Value = 1560
$ws.Cells.Item(98, 13).Value = -1775
$ws.Cells.Item(98, 14).Value = -4556
$ws.Cells.Item(113, 8).Value = 12346414
$ws.Cells.Item(113, 10).Value = 1103.5
$ws.Cells.Item(113, 12).Value = 3310.5
$ws.Cells.Item(113, 14).Value = -7650.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6884.6665
$ws.Cells.Item(70, 9).Value = 6923
$ws.Cells.Item(70, 10).Value = 6840.857
$ws.Cells.Item(70, 11).Value = 6923
$ws.Cells.Item(70, 12).Value = 6840.857
$ws.Cells.Item(70, 13).Value = -6653
$ws.Cells.Item(70, 14).Value = -7380.857
$ws.Cells.Item(73, 8).Value = 6884.6665
$ws.Cells.Item(73, 9).Value = 6923
$ws.Cells.Item(73, 10).Value = 6840.857
$ws.Cells.Item(73, 11).Value = 6923
$ws.Cells.Item(73, 12).Value = 6840.857
$ws.Cells.Item(73, 13).Value = -5987
$ws.Cells.Item(73, 14).Value = -8712.857
$ws.Cells.Item(113, 8).Value = 406760.28
$ws.Cells.Item(113, 9).Value = 596252.4399999999
$ws.Cells.Item(113, 10).Value = 4089.5
$ws.Cells.Item(113, 11).Value = 596252.4399999999
$ws.Cells.Item(113, 12).Value = 4089.5
$ws.Cells.Item(113, 13).Value = -594082.4399999999
$ws.Cells.Item(113, 14).Value = -8429.5
$ws.Cells.Item(121, 8).Value = 59657.8
$ws.Cells.Item(121, 10).Value = 59657.8
$ws.Cells.Item(121, 12).Value = 59657.8
$ws.Cells.Item(121, 14).Value = -63151.8
$ws.Cells.Item(132, 8).Value = 336084.88
$ws.Cells.Item(132, 10).Value = 46155.957
$ws.Cells.Item(132, 12).Value = 138467.871
$ws.Cells.Item(132, 14).Value = -143527.871

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 7039
$ws.Cells.Item(61, 9).Value = 4897.1113
$ws.Cells.Item(61, 11).Value = 4897.1113
$ws.Cells.Item(61, 13).Value = -4695.1113
$ws.Cells.Item(100, 8).Value = 76932
$ws.Cells.Item(100, 9).Value = 119390.22
$ws.Cells.Item(100, 11).Value = 119390.22
$ws.Cells.Item(100, 13).Value = -118849.22
$ws.Cells.Item(113, 8).Value = 7039
$ws.Cells.Item(113, 9).Value = 4897.1113
$ws.Cells.Item(113, 11).Value = 4897.1113
$ws.Cells.Item(113, 13).Value = -2727.1113
$ws.Cells.Item(122, 8).Value = 412984.84
$ws.Cells.Item(122, 9).Value = 5005.6
$ws.Cells.Item(122, 11).Value = 15016.8
$ws.Cells.Item(122, 13).Value = -12566.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 4020
$ws.Cells.Item(4, 9).Value = 5000
$ws.Cells.Item(4, 11).Value = 5000
$ws.Cells.Item(4, 13).Value = -4887

Write-Host "Applied all changes"